$d = $word.ActiveDocument

# Locate the paragraph that hosts the inline picture ("TS045_TC117 Step 1.png").
# Everything from that paragraph through the end of the document body
# (the picture paragraph, the page-break paragraph, and the "Passed"
# paragraph) is being removed, leaving the last table immediately
# followed by the section properties.
$startPos = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.InlineShapes.Count -gt 0) {
        $startPos = $p.Range.Start
        break
    }
}

if ($startPos -ge 0) {
    $endPos = $d.Content.End
    $rng = $d.Range($startPos, $endPos)
    $rng.Delete()
}
